# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The workbook's "K" column (column G) was recomputed from strikeouts (K)
# instead of the previous "Strike#" metric. This updates the already
# written values in column G for each data row to the newly computed K
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value
$kValues = [ordered]@{
    2  = 1
    3  = 1
    4  = 0
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 2
    29 = 1
    30 = 1
    31 = 1
    33 = 0
    34 = 0
    35 = 2
    36 = 2
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 1
    45 = 0
    46 = 1
    47 = 2
    49 = 1
    51 = 2
    52 = 1
    53 = 0
    54 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
